$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 161: correct the end time for the existing entry ---
$ws.Range("E161").Value = 0.4375

# --- Row 162: turn the previously-empty placeholder row into a full data row ---
$ws.Range("A162").Value = 2014
$ws.Range("B162").Value = 7
$ws.Range("C162").Value = 28
$ws.Range("D162").Value = 0.55555555555555558
$ws.Range("E162").Value = 0.57638888888888895
$ws.Range("F162").Formula = "=(E162-D162)*24*60"
$ws.Range("G162").Formula = "=F162/60"

# --- Row 163: new data row (start time only entered so far) ---
$ws.Range("A163").Value = 2014
$ws.Range("B163").Value = 7
$ws.Range("C163").Value = 28
$ws.Range("D163").NumberFormat = "hh:mm;@"
$ws.Range("D163").Value = 0.60416666666666663
$ws.Range("E163").Clear()
$ws.Range("E163").NumberFormat = "hh:mm;@"
$ws.Range("F163").Clear()
$ws.Range("F163").NumberFormat = "0"

# --- Row 164: blank spacer row (previously occupied by the summary rows) ---
$ws.Range("D164").NumberFormat = "hh:mm;@"
$ws.Range("E164").Clear()
$ws.Range("E164").NumberFormat = "hh:mm;@"
$ws.Range("F164").Clear()
$ws.Range("F164").NumberFormat = "0"

# --- Row 165-167: summary rows, shifted down by two rows, ranges updated ---
$ws.Range("E165").Value = "sum [min]"
$ws.Range("E165").HorizontalAlignment = -4152
$ws.Range("F165").Clear()
$ws.Range("F165").NumberFormat = "0"
$ws.Range("F165").Formula = "=SUM(F2:F162)"

$ws.Range("E166").Value = "sum [h]"
$ws.Range("E166").HorizontalAlignment = -4152
$ws.Range("F166").NumberFormat = "0.00"
$ws.Range("F166").Formula = "=F165/60"

$ws.Range("E167").Value = "sum [working weeks]"
$ws.Range("E167").HorizontalAlignment = -4152
$ws.Range("F167").NumberFormat = "0.00"
$ws.Range("F167").Formula = "=F166/38.5"

# --- Selection, matching the recorded cursor position after the edit ---
$ws.Range("E163").Select()

Write-Output "edit complete"
